# "export image feature implemented"
# Update the two data rows (row 2 and row 3) of the PT files export sheet
# with the new item data (SHIRT instead of TROUSER, new codes/sizes/prices).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a numeric-looking value while keeping it stored as TEXT
# (matches source data where these are shared-string values, not numbers).
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# ---------- Row 2 ----------
$ws.Range("A2").Value = "SHIRT"
$ws.Range("C2").Value = "S1"
$ws.Range("D2").Value = "A 1142/2"
$ws.Range("G2").Value = 4
$ws.Range("I2").ClearContents()
$ws.Range("J2").Value = "34/86 CM"
Set-TextValue $ws.Range("L2") "2222"
Set-TextValue $ws.Range("R2") "2699.0"
$ws.Range("T2").Value = 1

# ---------- Row 3 ----------
$ws.Range("A3").Value = "SHIRT"
$ws.Range("B3").Value = "C2"
$ws.Range("C3").Value = "S2"
Set-TextValue $ws.Range("D3") "3612"
$ws.Range("G3").Value = 5
$ws.Range("I3").Value = "OUTFIT BLACK"
Set-TextValue $ws.Range("J3") "42"
Set-TextValue $ws.Range("R3") "1111"
